$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.536.79"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "1.622.53"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.86"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.522"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.84%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.04"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.262"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0612"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0889"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").Value = "1.854.11"
$ws.Range("E12").Value = "  -0.64%  "
$ws.Range("D13").Value = "1.626.35"
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.04"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.549"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.42"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").Value = "27.538.19"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.10"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.54"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.66%  "
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.30"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.90"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("E24").Value = "  +6.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.15"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.87"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.23%  "
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.111"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.92%  "
$ws.Range("B28").Value = "BinanceUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.59"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("E30").Value = "  -0.37%  "
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("D33").Value = "1.457.33"
$ws.Range("E33").Value = "  +3.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.08"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.90%  "
$ws.Range("E35").Value = "  -2.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.34"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.949"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +6.73%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.565"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.89%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0168"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.35%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.865"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.37"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +7.10%  "
$ws.Range("E42").Value = "  +0.23%  "
$ws.Range("E43").Value = "  -1.75%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.23"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.93%  "
$ws.Range("D47").Value = "1.763.99"
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.18"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.16%  "
$ws.Range("D50").Value = "0.0₆0105"
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("E51").Value = "  +1.38%  "
